$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 100 (2020-02-27) and 101 (2020-02-28) appended to the data table.
# Columns B (date) and C (id) hold values that look numeric ("2020-02-27",
# "0213") but must stay text, matching the rest of the sheet -- force the
# "Text" number format while writing them, then restore the default style so
# the cells end up unstyled/text-typed just like their neighbours.
$textRange = $ws.Range("B100:C101")
$textRange.NumberFormat = "@"

# Row 100
$ws.Range("A100").Value = 1582761600
$ws.Range("B100").Value = "2020-02-27"
$ws.Range("C100").Value = "0213"
$ws.Range("D100").Value = "MTAG"
$ws.Range("E100").Value = 0.505
$ws.Range("F100").Value = 0.505
$ws.Range("G100").Value = 0.485
$ws.Range("H100").Value = 0.485
$ws.Range("I100").Value = 10414000

# Row 101
$ws.Range("A101").Value = 1582848000
$ws.Range("B101").Value = "2020-02-28"
$ws.Range("C101").Value = "0213"
$ws.Range("D101").Value = "MTAG"
$ws.Range("E101").Value = 0.47
$ws.Range("F101").Value = 0.475
$ws.Range("G101").Value = 0.44
$ws.Range("H101").Value = 0.45
$ws.Range("I101").Value = 10664700

$textRange.Style = "Normal"
